# Fix Training Data Issue
# The BF column holds a "Date" (as text) that was recorded using the
# wrong format (M-D-YYYY-YY, e.g. "6-2-2012-13"). Re-write it in
# ISO form (YYYY-MM-DD, e.g. "2013-06-02") for every data row.
#
# NumberFormat is forced to Text ("@") before the assignment so Excel's
# automatic date recognition doesn't silently turn the literal string
# into a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 31
$col      = 58   # column BF
$correctDate = "2013-06-02"

$rangeAddress = "BF" + $firstRow + ":BF" + $lastRow
$rng = $ws.Range($rangeAddress)
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $col).Value = $correctDate
}
